# Apply the cryptos-list price/volume refresh described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.486.19"
$ws.Range("E2").Value = "  -1.23%  "
$ws.Range("D3").Value = "2.347.44"
$ws.Range("E3").Value = "  +3.24%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("B5").Value = "XRP"
$ws.Range("C5").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D5").Value = "'0.652"
$ws.Range("E5").Value = "  +2.52%  "
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "'232.87"
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("D7").Value = "'66.10"
$ws.Range("E7").Value = "  +3.76%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.454"
$ws.Range("E9").Value = "  +1.37%  "
$ws.Range("D10").Value = "'0.0964"
$ws.Range("E10").Value = "  -2.86%  "
$ws.Range("E11").Value = "  -1.15%  "
$ws.Range("D12").Value = "'26.85"
$ws.Range("E12").Value = "  -2.06%  "
$ws.Range("D13").Value = "2.694.75"
$ws.Range("E13").Value = "  +3.20%  "
$ws.Range("D15").Value = "'15.48"
$ws.Range("E15").Value = "  -1.73%  "
$ws.Range("D16").Value = "'6.28"
$ws.Range("E16").Value = "  +2.55%  "
$ws.Range("D17").Value = "'0.852"
$ws.Range("E17").Value = "  +1.70%  "
$ws.Range("D18").Value = "2.340.68"
$ws.Range("E18").Value = "  +3.66%  "
$ws.Range("D19").Value = "43.428.47"
$ws.Range("E19").Value = "  -0.98%  "
$ws.Range("D20").Value = "0.0₃0985"
$ws.Range("E20").Value = "  -1.71%  "
$ws.Range("D21").Value = "'74.35"
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("D22").Value = "'6.24"
$ws.Range("E22").Value = "  +2.33%  "
$ws.Range("D23").Value = "'249.74"
$ws.Range("E23").Value = "  -1.09%  "
$ws.Range("E24").Value = "  +16.21%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  -0.80%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "'2.28"
$ws.Range("E27").Value = "  +0.82%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "'9.95"
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("D29").Value = "'175.05"
$ws.Range("E29").Value = "  +1.90%  "
$ws.Range("D30").Value = "'22.21"
$ws.Range("E30").Value = "  +6.18%  "
$ws.Range("E31").Value = "  +6.69%  "
$ws.Range("E32").Value = "  -6.58%  "
$ws.Range("E33").Value = "  +0.78%  "
$ws.Range("D34").Value = "'4.99"
$ws.Range("E34").Value = "  +3.87%  "
$ws.Range("D35").Value = "'0.0689"
$ws.Range("E35").Value = "  -1.10%  "
$ws.Range("D36").Value = "'4.97"
$ws.Range("E36").Value = "  +1.64%  "
$ws.Range("D37").Value = "'2.53"
$ws.Range("E37").Value = "  +9.33%  "
$ws.Range("E38").Value = "  -0.66%  "
$ws.Range("E39").Value = "  -4.91%  "
$ws.Range("E40").Value = "  -2.31%  "
$ws.Range("E41").Value = "  +8.87%  "
$ws.Range("D43").Value = "'18.10"
$ws.Range("E43").Value = "  +2.57%  "
$ws.Range("E44").Value = "  +7.95%  "
$ws.Range("D45").Value = "'99.44"
$ws.Range("E45").Value = "  +1.25%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").Value = "'0.0948"
$ws.Range("E47").Value = "  -2.92%  "
$ws.Range("E48").Value = "  -0.69%  "
$ws.Range("D49").Value = "1.440.33"
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("E50").Value = "  -10.71%  "
$ws.Range("D51").Value = "'9.85"
$ws.Range("E51").Value = "  -5.43%  "
